$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "SCD0020"

# Update the TC_ID value
$ws.Range("B2").Value = "SCD0020-005"

# Apply left/center alignment to every used cell in the sheet
$used = $ws.UsedRange
$used.HorizontalAlignment = -4131   # xlLeft
$used.VerticalAlignment = -4108     # xlCenter
